# Applies the cryptos-list price/volume refresh described in the commit message
# ("Updated cryptos list on ... with GitHub Actions"): refreshed Price (D) and
# Volume(1h) (E) figures for each coin row, plus a rank swap between the VeChain
# and Maker rows (36/37) where name, link, price and volume all move together.
#
# Every touched cell is plain text in the workbook (Price column includes values
# like "26.396.99" / "0.0821" that are not valid numbers, or that Excel would
# re-type and strip trailing zeros from, e.g. "0.250" -> 0.25). Numeric-looking
# Price values are therefore written with a leading apostrophe to force Excel to
# keep them as text, then the style is reset to Normal so the one-off "quote
# prefix" formatting that the apostrophe trick applies does not stick to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.329.48'
$ws.Range('E2').Value = '  +4.00%  '
$ws.Range('D3').Value = '1.608.46'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '''213.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.76%  '
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  +2.13%  '
$ws.Range('D8').Value = '''0.250'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.94%  '
$ws.Range('D9').Value = '''0.0620'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('D10').Value = '''18.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('D11').Value = '''0.0820'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.82%  '
$ws.Range('D12').Value = '1.835.68'
$ws.Range('E12').Value = '  +2.73%  '
$ws.Range('D13').Value = '1.611.31'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('D14').Value = '''4.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '''0.514'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').Value = '26.331.90'
$ws.Range('E16').Value = '  +4.05%  '
$ws.Range('D17').Value = '''60.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  +2.59%  '
$ws.Range('D19').Value = '''209.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.85%  '
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').Value = '''4.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('D22').Value = '''9.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('D23').Value = '''6.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.63%  '
$ws.Range('D24').Value = '''1.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.64%  '
$ws.Range('D25').Value = '''142.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('E27').Value = '  -3.78%  '
$ws.Range('D28').Value = '''15.31'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.96%  '
$ws.Range('D29').Value = '''6.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('D31').Value = '''0.0473'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('E33').Value = '  +0.89%  '
$ws.Range('D34').Value = '''1.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').Value = '''2.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.110.55'
$ws.Range('E36').Value = '  +1.85%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.0161'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.19%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').Value = '''0.788'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.69%  '
$ws.Range('D41').Value = '''0.499'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').Value = '''0.779'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = '1.748.57'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').Value = '''93.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E46').Value = '  +10.15%  '
$ws.Range('E47').Value = '  -5.09%  '
$ws.Range('D48').Value = '''53.70'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('E51').Value = '  -0.28%  '
